# Weekly fruit/vegetable price update.
# Insert 4 new rows (new week of data, 2021-11-23 = serial 44523) right
# before the current row 93, pushing the existing rows 93:107 down to
# 97:111.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 93:107 down by 4 rows to make room for the new
# weekly entries.
$ws.Range("A93:R96").EntireRow.Insert()

# Shared / constant values for this product block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 300000000
$categoria = "Espárragos"
$clasif    = "Hortaliza"

$fecha = 44523

$rows = @(
    @{ Row=93; Calidad="Banquete"; Vol=250; PMin=1600; PMax=1600; PProm=1600 },
    @{ Row=94; Calidad="Primera";  Vol=340; PMin=1400; PMax=1400; PProm=1400 },
    @{ Row=95; Calidad="Segunda";  Vol=160; PMin=1200; PMax=1200; PProm=1200 },
    @{ Row=96; Calidad="Tercera";  Vol=106; PMin=1000; PMax=1000; PProm=1000 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $catId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = "Sin especificar"
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Vol
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = "`$/kilo"
    $ws.Cells.Item($row, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = $clasif
}
